$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F53").Value = "N. crash:"
$ws.Range("G53").Formula = '=COUNTIF(G2:G51,"True")'

$ws.Range("F54").Value = "Mean:"
$ws.Range("H54").Formula = "=AVERAGE(H2:H51)"
$ws.Range("I54:P54").Formula = "=AVERAGE(I2:I51)"

$ws.Range("F55").Value = "Standard Deviation:"
$ws.Range("H55").Formula = "=STDEV.S(H2:H51)"
$ws.Range("I55:P55").Formula = "=STDEV.S(I2:I51)"

$ws.Range("T52").Select() | Out-Null
